# Scheduled runner update: refresh cached market-price / profit figures
# (currentAveragePrice*, Leve*Price*, Leve*Profit*) across several leve
# tracker sheets. Values mirror the latest Universalis price snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2613.1738
$ws.Range("I40").Value = 2533.5715
$ws.Range("K40").Value = 2533.5715
$ws.Range("M40").Value = -2358.5715

$ws.Range("H62").Value = 3528.75
$ws.Range("I62").Value = 2511.1765
$ws.Range("J62").Value = 6000
$ws.Range("K62").Value = 2511.1765
$ws.Range("L62").Value = 6000
$ws.Range("M62").Value = -1887.1765
$ws.Range("N62").Value = -7248

$ws.Range("H65").Value = 3528.75
$ws.Range("I65").Value = 2511.1765
$ws.Range("J65").Value = 6000
$ws.Range("K65").Value = 12555.8825
$ws.Range("L65").Value = 30000
$ws.Range("M65").Value = -9435.8825
$ws.Range("N65").Value = -36240

$ws.Range("H132").Value = 35870.453
$ws.Range("I132").Value = 63940.234
$ws.Range("J132").Value = 1785.7142
$ws.Range("K132").Value = 191820.702
$ws.Range("L132").Value = 5357.142599999999
$ws.Range("M132").Value = -189290.702
$ws.Range("N132").Value = -10417.1426

$ws.Range("H138").Value = 3801.04
$ws.Range("I138").Value = 2329.6667
$ws.Range("J138").Value = 4628.6875
$ws.Range("K138").Value = 6989.000100000001
$ws.Range("L138").Value = 13886.0625
$ws.Range("M138").Value = -1849.000100000001
$ws.Range("N138").Value = -24166.0625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1926.1666
$ws.Range("I88").Value = 1300
$ws.Range("J88").Value = 2239.25
$ws.Range("K88").Value = 1300
$ws.Range("L88").Value = 2239.25
$ws.Range("M88").Value = -894
$ws.Range("N88").Value = -3051.25

$ws.Range("H91").Value = 1926.1666
$ws.Range("I91").Value = 1300
$ws.Range("J91").Value = 2239.25
$ws.Range("K91").Value = 1300
$ws.Range("L91").Value = 2239.25
$ws.Range("M91").Value = 104
$ws.Range("N91").Value = -5047.25

$ws.Range("H132").Value = 1577.3818
$ws.Range("I132").Value = 1339.102
$ws.Range("J132").Value = 3523.3333
$ws.Range("K132").Value = 4017.306
$ws.Range("L132").Value = 10569.9999
$ws.Range("M132").Value = -1487.306
$ws.Range("N132").Value = -15629.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 497.48
$ws.Range("I94").Value = 479.83334
$ws.Range("J94").Value = 542.8570999999999
$ws.Range("K94").Value = 479.83334
$ws.Range("L94").Value = 542.8570999999999
$ws.Range("M94").Value = -28.83334000000002
$ws.Range("N94").Value = -1444.8571

$ws.Range("H99").Value = 1566.3334
$ws.Range("I99").Value = 1197.4
$ws.Range("J99").Value = 2488.6667
$ws.Range("K99").Value = 1197.4
$ws.Range("L99").Value = 2488.6667
$ws.Range("M99").Value = 300.5999999999999
$ws.Range("N99").Value = -5484.6667

$ws.Range("H105").Value = 2198.2354
$ws.Range("I105").Value = 1922.8572
$ws.Range("J105").Value = 3483.3333
$ws.Range("K105").Value = 1922.8572
$ws.Range("L105").Value = 3483.3333
$ws.Range("M105").Value = -175.8571999999999
$ws.Range("N105").Value = -6977.3333

$ws.Range("H134").Value = 1563.5758
$ws.Range("I134").Value = 1342.3334
$ws.Range("J134").Value = 1950.75
$ws.Range("K134").Value = 4027.0002
$ws.Range("L134").Value = 5852.25
$ws.Range("M134").Value = -1492.0002
$ws.Range("N134").Value = -10922.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 200
$ws.Range("I2").Value = 200
$ws.Range("K2").Value = 200
$ws.Range("M2").Value = -87

$ws.Range("H15").Value = 300
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

$ws.Range("H31").Value = 2915.99
$ws.Range("I31").Value = 1762.4576
$ws.Range("J31").Value = 4575.951
$ws.Range("K31").Value = 1762.4576
$ws.Range("L31").Value = 4575.951
$ws.Range("M31").Value = -1467.4576
$ws.Range("N31").Value = -5165.951

$ws.Range("H34").Value = 2915.99
$ws.Range("I34").Value = 1762.4576
$ws.Range("J34").Value = 4575.951
$ws.Range("K34").Value = 1762.4576
$ws.Range("L34").Value = 4575.951
$ws.Range("M34").Value = -1560.4576
$ws.Range("N34").Value = -4979.951

$ws.Range("H62").Value = 3880.7144
$ws.Range("I62").Value = 3135.9092
$ws.Range("J62").Value = 4700
$ws.Range("K62").Value = 3135.9092
$ws.Range("L62").Value = 4700
$ws.Range("M62").Value = -2511.9092
$ws.Range("N62").Value = -5948

$ws.Range("H65").Value = 3880.7144
$ws.Range("I65").Value = 3135.9092
$ws.Range("J65").Value = 4700
$ws.Range("K65").Value = 15679.546
$ws.Range("L65").Value = 23500
$ws.Range("M65").Value = -12559.546
$ws.Range("N65").Value = -29740

$ws.Range("H105").Value = 1921.5385
$ws.Range("I105").Value = 2525
$ws.Range("J105").Value = 956
$ws.Range("K105").Value = 2525
$ws.Range("L105").Value = 956
$ws.Range("M105").Value = -778
$ws.Range("N105").Value = -4450

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 1250.4117
$ws.Range("I125").Value = 609.5
$ws.Range("J125").Value = 1600
$ws.Range("K125").Value = 1828.5
$ws.Range("L125").Value = 4800
$ws.Range("M125").Value = 3091.5
$ws.Range("N125").Value = -14640

$ws.Range("H131").Value = 760.72
$ws.Range("J131").Value = 842.84705
$ws.Range("L131").Value = 2528.54115
$ws.Range("N131").Value = -12608.54115

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 50000
$ws.Range("J4").Value = 50000
$ws.Range("L4").Value = 50000
$ws.Range("N4").Value = -50224

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 50948.145
$ws.Range("I7").Value = 74321.86
$ws.Range("K7").Value = 74321.86
$ws.Range("M7").Value = -74209.86

$ws.Range("H14").Value = 23263.334
$ws.Range("J14").Value = 23263.334
$ws.Range("L14").Value = 23263.334
$ws.Range("N14").Value = -23607.334

$ws.Range("H40").Value = 36792.785
$ws.Range("I40").Value = 39553.77
$ws.Range("J40").Value = 900
$ws.Range("K40").Value = 39553.77
$ws.Range("L40").Value = 900
$ws.Range("M40").Value = -39417.77
$ws.Range("N40").Value = -1172

$ws.Range("H68").Value = 1629.4445
$ws.Range("I68").Value = 1164
$ws.Range("J68").Value = 1862.1666
$ws.Range("K68").Value = 1164
$ws.Range("L68").Value = 1862.1666
$ws.Range("M68").Value = -415
$ws.Range("N68").Value = -3360.1666

$ws.Range("H71").Value = 1629.4445
$ws.Range("I71").Value = 1164
$ws.Range("J71").Value = 1862.1666
$ws.Range("K71").Value = 5820
$ws.Range("L71").Value = 9310.833000000001
$ws.Range("M71").Value = -2076
$ws.Range("N71").Value = -16798.833

$ws.Range("H122").Value = 22224678
$ws.Range("I122").Value = 55556610
$ws.Range("J122").Value = 3393.3333
$ws.Range("K122").Value = 166669830
$ws.Range("L122").Value = 10179.9999
$ws.Range("M122").Value = -166667380
$ws.Range("N122").Value = -15079.9999

$ws.Range("H126").Value = 50948.145
$ws.Range("I126").Value = 74321.86
$ws.Range("K126").Value = 222965.58
$ws.Range("M126").Value = -220495.58

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3300
$ws.Range("I62").Value = 3483.3333
$ws.Range("J62").Value = 2933.3333
$ws.Range("K62").Value = 3483.3333
$ws.Range("L62").Value = 2933.3333
$ws.Range("M62").Value = -2859.3333
$ws.Range("N62").Value = -4181.3333

$ws.Range("H64").Value = 48730
$ws.Range("J64").Value = 48730
$ws.Range("L64").Value = 48730
$ws.Range("N64").Value = -49226

$ws.Range("H65").Value = 3300
$ws.Range("I65").Value = 3483.3333
$ws.Range("J65").Value = 2933.3333
$ws.Range("K65").Value = 17416.6665
$ws.Range("L65").Value = 14666.6665
$ws.Range("M65").Value = -14296.6665
$ws.Range("N65").Value = -20906.6665

$ws.Range("H67").Value = 48730
$ws.Range("J67").Value = 48730
$ws.Range("L67").Value = 48730
$ws.Range("N67").Value = -50446

$ws.Range("H122").Value = 79369.234
$ws.Range("I122").Value = 168833.33
$ws.Range("J122").Value = 2685.7144
$ws.Range("K122").Value = 506499.99
$ws.Range("L122").Value = 8057.1432
$ws.Range("M122").Value = -504049.99
$ws.Range("N122").Value = -12957.1432

$ws.Range("H136").Value = 3085.746
$ws.Range("I136").Value = 621.27026
$ws.Range("J136").Value = 6592.885
$ws.Range("K136").Value = 1863.81078
$ws.Range("L136").Value = 19778.655
$ws.Range("M136").Value = 686.18922
$ws.Range("N136").Value = -24878.655
